# "Minor edit to last slide" -- also carries the small Bootstrap CSS
# bullet tweak on slide 2 that shipped in the same save.
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 2 ("Technologies Used" style slide): bullet 3 gains " forms"
#   "Bootstrap CSS – for styling"  ->  "Bootstrap CSS – for styling forms"
# ---------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$shape2 = $slide2.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange
$para2 = $tr2.Paragraphs(3, 1)
$offset2 = $para2.Text.IndexOf("styling")
$word2 = $tr2.Characters($para2.Start + $offset2, 7)
$word2.Text = "styling forms"

# ---------------------------------------------------------------
# Slide 9 (last slide, "Individual Contributions"): rewrite Mikah's
# bullet with the updated wording.
#   "Mikah: Got us converted from basic CSS to Bootstrap, did the
#    other few of the original pages shown at the midterm
#    presentation."
#   ->
#   "Mikah: Applied Bootstrap to the form components, did the other
#    few of the original webpages shown at the midterm presentation.
#    Some front-end and the breadcrumbs."
# ---------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$shape9 = $slide9.Shapes.Item(2)
$tr9 = $shape9.TextFrame.TextRange
$para9 = $tr9.Paragraphs(3, 1)
$base9 = $para9.Start

# Replace piece by piece, right-to-left, so earlier offsets stay valid.
$seg6 = $tr9.Characters($base9 + $para9.Text.IndexOf("."), 1)
$seg6.Text = ". Some front-end and the breadcrumbs."

$para9 = $tr9.Paragraphs(3, 1)
$seg5Start = $para9.Text.IndexOf("shown at the midterm presentation")
$seg5 = $tr9.Characters($base9 + $seg5Start, "shown at the midterm presentation".Length)
$seg5.Text = "shown at the midterm presentation"

$para9 = $tr9.Paragraphs(3, 1)
$seg4Start = $para9.Text.IndexOf("pages ")
$seg4 = $tr9.Characters($base9 + $seg4Start, "pages ".Length)
$seg4.Text = "webpages "

$para9 = $tr9.Paragraphs(3, 1)
$seg3Start = $para9.Text.IndexOf("did the other few of the original ")
$seg3 = $tr9.Characters($base9 + $seg3Start, "did the other few of the original ".Length)
$seg3.Text = "did the other few of the original "

$para9 = $tr9.Paragraphs(3, 1)
$seg2Start = $para9.Text.IndexOf("Got us converted from basic CSS to Bootstrap, ")
$seg2 = $tr9.Characters($base9 + $seg2Start, "Got us converted from basic CSS to Bootstrap, ".Length)
$seg2.Text = "Applied Bootstrap to the form components, "

$para9 = $tr9.Paragraphs(3, 1)
$seg1 = $tr9.Characters($base9 + 0, "Mikah: ".Length)
$seg1.Text = "Mikah: "
